$d = $word.ActiveDocument

function Split-Run($start, $end) {
    $r = $d.Range($start, $end)
    $r.Font.Bold = 1
    $r.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "...alpha is multiplied by in stochastic GD is exactly..."
#    -> split "alpha is multiplied by in stochastic GD " into three runs,
#       capitalising "Stochastic".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("alpha is multiplied by in stochastic GD ") | Out-Null
$s = $rng.Start
$aLen = "alpha is multiplied by in ".Length
$bLen = "Stochastic ".Length
$p2Start = $s + $aLen
$p2End = $p2Start + $bLen
$part2 = $d.Range($p2Start, $p2End)
$part2.Text = "Stochastic "
Split-Run $p2Start $p2End

Write-Output "1 ok"

# ---------------------------------------------------------------------
# 2) "What Stochastic GD is doing is scan through..."
#    -> drop the leading "What " run, and turn "is doing is scan" into
#       two runs: "=" and " scan".
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("What Stochastic GD is doing is scan") | Out-Null
$s = $rng.Start
$whatLen = "What ".Length
$delRange = $d.Range($s, $s + $whatLen)
$delRange.Text = ""

$stoLen = "Stochastic ".Length
$gdLen = "GD ".Length
$stoEnd = $s + $stoLen
Split-Run $s $stoEnd
$gdStart = $stoEnd
$gdEnd = $gdStart + $gdLen
Split-Run $gdStart $gdEnd

$oldLen = "is doing is scan".Length
$scanRangeEnd = $gdEnd + $oldLen
$scanRange = $d.Range($gdEnd, $scanRangeEnd)
$scanRange.Text = "= scan"
$eqEnd = $gdEnd + 1
Split-Run $gdEnd $eqEnd

Write-Output "2 ok"

# ---------------------------------------------------------------------
# 3) "...training example {x(1), y(1)}..." -> "...training {x(1), y(1)}..."
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("the 1st training example {x(1), y(1)}") | Out-Null
$s = $rng.Start
$prefixLen = "the 1st ".Length
$oldLen = "training example ".Length
$replStart = $s + $prefixLen
$replEnd = $replStart + $oldLen
$replRange = $d.Range($replStart, $replEnd)
Write-Output "P3 before: [$($replRange.Text)]"
$replRange.Text = "training "

Write-Output "3 ok"

# ---------------------------------------------------------------------
# 4) "example, take a small gradient descent step" ->
#    "example, take a <i>small</i> gradient descent step" (space split off
#    into its own, non-italic run).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("example, take a small gradient descent step") | Out-Null
$s = $rng.Start
$prefixLen = "example, take a ".Length
$smallLen = "small".Length
$smallStart = $s + $prefixLen
$smallEnd = $smallStart + $smallLen
$smallRange = $d.Range($smallStart, $smallEnd)
Write-Output "P4 before: [$($smallRange.Text)]"
$smallRange.Font.Italic = 1

Write-Output "4 ok"

# ---------------------------------------------------------------------
# 5) italic "training example" followed by ". " -> ". " loses its
#    trailing space.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("training example. ") | Out-Null
$s = $rng.Start
$prefixLen = "training example".Length
$periodStart = $s + $prefixLen
$periodEnd = $periodStart + ". ".Length
$periodRange = $d.Range($periodStart, $periodEnd)
Write-Output "P5 before: [$($periodRange.Text)]"
$periodRange.Text = "."

Write-Output "5 ok"

# ---------------------------------------------------------------------
# 6) "In other words, it looks at the 1st example" ->
#    "In other words: look at the 1st example"
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("other words, it looks at the") | Out-Null
$s = $rng.Start
$oldLen = "other words, it looks ".Length
$wordsRange = $d.Range($s, $s + $oldLen)
Write-Output "P6 before: [$($wordsRange.Text)]"
$wordsRange.Text = "other words: look "

Write-Output "6 ok"
